# "adding manual review notes from 20200104"
# Row 19 (Brent_2b_... fastq file) gets a manual-review audit flag and note.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# manualAudit (H19): 0 -> 1
$ws.Range("H19").Value = 1

# manualStatus (I19): new note "[512]"
$ws.Range("I19").Value = "[512]"

# Move the selection to the newly-annotated cell, matching the reviewer's
# last position when they saved the file.
$ws.Range("I19").Select()
